# Generate Report for Handback
# Update status text and timestamps to reflect a new handback report run,
# and widen the two "in sync" status columns to fit the longer text.

$wb = $excel.ActiveWorkbook

$wsOverview = $wb.Worksheets.Item("Overview")
$wsZhCn = $wb.Worksheets.Item("zh-cn")
$wsDeDe = $wb.Worksheets.Item("de-de")

# --- Overview sheet ---
# Columns E and F hold the "Handed back: in sync with en-US" status, column G the timestamp.
$wsOverview.Range("E2").Value = "Handed back: not in sync with en-US"
$wsOverview.Range("F2").Value = "Handed back: not in sync with en-US"
$wsOverview.Range("G2").Value = "2017-02-10 07:32:47"

# Widen columns E and F to fit the longer status text.
$wsOverview.Columns.Item(5).ColumnWidth = 33.4602203369141
$wsOverview.Columns.Item(6).ColumnWidth = 33.4602203369141

# --- zh-cn sheet ---
$wsZhCn.Range("C2").Value = "Handed back: not in sync with en-US"
$wsZhCn.Range("H2").Value = "2017-02-10 07:32:30"
$wsZhCn.Range("L2").Value = "2017-02-10 07:34:12"
$wsZhCn.Columns.Item(3).ColumnWidth = 33.4602203369141

# --- de-de sheet ---
$wsDeDe.Range("C2").Value = "Handed back: not in sync with en-US"
$wsDeDe.Range("H2").Value = "2017-02-10 07:32:47"
$wsDeDe.Range("L2").Value = "2017-02-10 07:34:34"
$wsDeDe.Columns.Item(3).ColumnWidth = 33.4602203369141
